# Insert a new price record for "Femacal de La Calera" (Frutilla) as row 82,
# shifting the existing rows 82-114 down to 83-115 (dimension grows to A1:T115).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by inserting a fresh row at position 82.
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(82, 1).Value  = 3
$ws.Cells.Item(82, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(82, 3).Value  = "Coquimbo"
$ws.Cells.Item(82, 4).Value  = 44468
$ws.Cells.Item(82, 5).Value  = 5
$ws.Cells.Item(82, 6).Value  = "Fruta"
$ws.Cells.Item(82, 7).Value  = 100101
$ws.Cells.Item(82, 8).Value  = "Berries"
$ws.Cells.Item(82, 9).Value  = 100112025
$ws.Cells.Item(82, 10).Value = "Frutilla"
$ws.Cells.Item(82, 11).Value = "Sin especificar"
$ws.Cells.Item(82, 12).Value = "Especial"
$ws.Cells.Item(82, 13).Value = 65
$ws.Cells.Item(82, 14).Value = 15000
$ws.Cells.Item(82, 15).Value = 15000
$ws.Cells.Item(82, 16).Value = 15000
$ws.Cells.Item(82, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(82, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(82, 19).Value = 2143
$ws.Cells.Item(82, 20).Value = 7
